$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Typist / Typist QC) before the old "Client" column (E),
# shifting Client..Tier from E:M to G:O.
$ws.Columns("E:F").Insert()

# New header cells
$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"

# Row 2 (order FTC18-001) - Typist / Typist QC values
$ws.Range("E2").Value = "SIPL5317"
$ws.Range("F2").Value = "SIPL5317"
# Old "Search(T1)" value (now shifted to O2) is no longer used
$ws.Range("O2").Value = ""

# Row 3 (order FTC18-002) - employee ids moved out of C3/D3, Typist / Typist QC added
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"
$ws.Range("N3").Value = "Typing"
$ws.Range("O3").Value = "Typing(T1)"

# Move the visible selection like the saved workbook
$ws.Range("J7").Select()
